$d = $word.ActiveDocument

$replacements = @(
    @{old="165÷5="; new="259÷7="},
    @{old="120÷3="; new="275÷7="},
    @{old="529÷7="; new="443÷2="},
    @{old="782÷6="; new="178÷7="},
    @{old="538÷7="; new="126÷4="},
    @{old="178÷2="; new="568÷4="},
    @{old="185÷5="; new="193÷8="},
    @{old="883÷6="; new="515÷8="},
    @{old="131÷3="; new="501÷8="},
    @{old="506÷2="; new="286÷5="},
    @{old="117÷6="; new="413÷5="},
    @{old="991÷2="; new="518÷8="},
    @{old="936÷4="; new="279÷6="},
    @{old="161÷4="; new="905÷3="},
    @{old="728÷8="; new="191÷6="},
    @{old="952÷4="; new="825÷4="},
    @{old="787÷3="; new="132÷3="},
    @{old="673÷7="; new="927÷9="},
    @{old="435÷8="; new="340÷9="},
    @{old="838÷7="; new="287÷7="},
    @{old="150÷5="; new="741÷5="},
    @{old="802÷4="; new="320÷6="},
    @{old="741÷4="; new="614÷6="},
    @{old="110÷5="; new="827÷8="},
    @{old="194÷4="; new="336÷7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
